$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns D that receive new text values need to stay text (not auto-converted to numbers).
# Pre-format the touched D cells as Text, then restore the original (default/"Normal") style
# after writing, so no residual numeric formatting is left behind.
$dCells = @("D2","D3","D4","D5","D6","D11","D13","D14","D16","D17","D19","D22","D23","D24","D27","D29","D31","D35","D36","D39","D40","D42","D45","D48","D50")
foreach ($addr in $dCells) { $ws.Range($addr).NumberFormat = "@" }

$ws.Range("D2").Value = '70.164.39'
$ws.Range("E2").Value = '  +0.42%  '
$ws.Range("D3").Value = '3.603.13'
$ws.Range("E3").Value = '  +2.27%  '
$ws.Range("D4").Value = '0.999'
$ws.Range("E4").Value = '  +0.16%  '
$ws.Range("D5").Value = '604.80'
$ws.Range("E5").Value = '  +0.46%  '
$ws.Range("D6").Value = '196.21'
$ws.Range("E6").Value = '  +0.37%  '
$ws.Range("E7").Value = '  +0.12%  '
$ws.Range("E8").Value = '  +0.07%  '
$ws.Range("E9").Value = '  -2.17%  '
$ws.Range("E10").Value = '  -0.95%  '
$ws.Range("D11").Value = '53.77'
$ws.Range("E11").Value = '  -0.49%  '
$ws.Range("E12").Value = '  +0.46%  '
$ws.Range("D13").Value = '9.57'
$ws.Range("E13").Value = '  +0.39%  '
$ws.Range("D14").Value = '4.178.45'
$ws.Range("E14").Value = '  +2.48%  '
$ws.Range("E15").Value = '  +3.71%  '
$ws.Range("D16").Value = '594.70'
$ws.Range("E16").Value = '  -1.35%  '
$ws.Range("D17").Value = '70.335.52'
$ws.Range("E17").Value = '  +0.42%  '
$ws.Range("E18").Value = '  -0.10%  '
$ws.Range("D19").Value = '3.607.45'
$ws.Range("E19").Value = '  +2.71%  '
$ws.Range("E20").Value = '  +1.32%  '
$ws.Range("E21").Value = '  +0.06%  '
$ws.Range("D22").Value = '17.75'
$ws.Range("E22").Value = '  -3.04%  '
$ws.Range("D23").Value = '5.18'
$ws.Range("E23").Value = '  -1.06%  '
$ws.Range("D24").Value = '101.92'
$ws.Range("E25").Value = '  +0.11%  '
$ws.Range("E26").Value = '  -1.63%  '
$ws.Range("D27").Value = '10.73'
$ws.Range("E27").Value = '  -1.84%  '
$ws.Range("E28").Value = '  -0.90%  '
$ws.Range("D29").Value = '33.78'
$ws.Range("E29").Value = '  +0.45%  '
$ws.Range("E30").Value = '  +5.07%  '
$ws.Range("D31").Value = '7.14'
$ws.Range("E31").Value = '  +0.31%  '
$ws.Range("E32").Value = '  -3.43%  '
$ws.Range("E33").Value = '  +0.78%  '
$ws.Range("E34").Value = '  +0.12%  '
$ws.Range("D35").Value = '0.0₃0888'
$ws.Range("E35").Value = '  +7.26%  '
$ws.Range("D36").Value = '3.913.91'
$ws.Range("E36").Value = '  +3.88%  '
$ws.Range("E37").Value = '  +0.52%  '
$ws.Range("E38").Value = '  +0.04%  '
$ws.Range("D39").Value = '516.84'
$ws.Range("E39").Value = '  +5.68%  '
$ws.Range("D40").Value = '36.86'
$ws.Range("E40").Value = '  +0.25%  '
$ws.Range("E41").Value = '  -0.93%  '
$ws.Range("D42").Value = '3.54'
$ws.Range("E42").Value = '  -1.34%  '
$ws.Range("E43").Value = '  -2.15%  '
$ws.Range("E44").Value = '  -0.54%  '
$ws.Range("D45").Value = '3.41'
$ws.Range("E45").Value = '  +2.32%  '
$ws.Range("E46").Value = '  +1.15%  '
$ws.Range("E47").Value = '  +0.01%  '
$ws.Range("D48").Value = '8.62'
$ws.Range("E48").Value = '  -0.45%  '
$ws.Range("E49").Value = '  -0.18%  '
$ws.Range("D50").Value = '0.000251'
$ws.Range("E50").Value = '  +2.70%  '
$ws.Range("E51").Value = '  -2.25%  '

foreach ($addr in $dCells) { $ws.Range($addr).Style = "Normal" }
